$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new reservation row (row 31): activity description + hours
$ws.Range("A31").Value = "4. iterace - další implementace (vybrat stůl)"
$ws.Range("A31").Style = $ws.Range("A30").Style
$ws.Range("B31").Value = 1

# Move active selection to A32 (next empty row), as in the author's edit
$ws.Range("A32").Select()
